$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "Gino"
$ws.Range("B6").Value = "Gino"
$ws.Range("B7").Value = "Gino "
$ws.Range("B8").Value = "Gino"

$ws.Range("C5").Value = "NO"
$ws.Range("C6").Value = "NO"
$ws.Range("C7").Value = "YES"
$ws.Range("C8").Value = "NO"

$ws.Range("C5").Select()
